# Adds a new "MGI Mammalian Phenotype Level 4:" (Heading 3) section with a
# follow-up paragraph of body text, inserted right after the
# "Neurotrophin signaling pathway" paragraph (consuming one of the blank
# paragraphs that follow it). The trailing "_GoBack" bookmark (which Word
# keeps at the most-recently-edited spot) is moved onto the new text by
# simply re-adding a bookmark with that name at the new location; Word
# bookmark names are unique, so this implicitly removes the old
# "_GoBack" bookmark that used to sit on "Merge PTM and Gene Expression
# Data".

$d = $word.ActiveDocument

# Locate the (currently empty) paragraph that is the second blank
# paragraph following "Neurotrophin signaling pathway" -- this is the
# one that gets turned into the new two-paragraph section.
$anchor = $d.Content
$anchor.Find.Execute("Neurotrophin signaling pathway", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara = $anchor.Paragraphs.Item(1)
$targetIndex = $anchorPara.Index + 2

$targetPara = $d.Paragraphs.Item($targetIndex)
$targetStart = $targetPara.Range.Start

# Insert a brand-new paragraph ("MGI Mammalian Phenotype Level 4:") right
# before the blank paragraph; InsertXML creates a proper new paragraph
# mark so the old blank paragraph is pushed one slot later instead of
# being overwritten.
$insertPoint = $d.Range($targetStart, $targetStart)
$headingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>MGI Mammalian Phenotype Level 4:</w:t></w:r></w:p>'
$insertPoint.InsertXML($headingXml)

$headingPara = $d.Paragraphs.Item($targetIndex)
$headingPara.Style = "Heading 3"

# The paragraph right after the new heading is the original blank
# paragraph; give it the body text.
$bodyPara = $d.Paragraphs.Item($targetIndex + 1)
$bodyRange = $bodyPara.Range
$placeholder = "Abnormal embryo for both clusters. Z"
$bodyRange.Text = $placeholder

# Re-fetch the paragraph range (text was replaced) and drop a bookmark
# immediately before the trailing placeholder character -- placing a
# collapsed bookmark exactly at a paragraph's final text boundary isn't
# reliable, so we bookmark next to a throwaway character and then erase
# that character, leaving the bookmark sitting right after "clusters. ".
$bodyPara2 = $d.Paragraphs.Item($targetIndex + 1)
$bodyRange2 = $bodyPara2.Range
$bmPos = $bodyRange2.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholderRange = $d.Range($bodyRange2.End - 2, $bodyRange2.End - 1)
$placeholderRange.Text = ""
